# Update the grouped_recurrence_score filter value from "16-20" to "31-35"
# in both Cypher queries stored on the "startup" sheet (columns B and C,
# rows 2-4 all share the same two query strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("startup")

$newQueryB = $ws.Range("B2").Value2.Replace('sf.grouped_recurrence_score IN ["16-20"]', 'sf.grouped_recurrence_score IN ["31-35"]')
$newQueryC = $ws.Range("C2").Value2.Replace('sf.grouped_recurrence_score IN ["16-20"]', 'sf.grouped_recurrence_score IN ["31-35"]')

$ws.Range("B2").Value2 = $newQueryB
$ws.Range("B3").Value2 = $newQueryB
$ws.Range("B4").Value2 = $newQueryB

$ws.Range("C2").Value2 = $newQueryC
$ws.Range("C3").Value2 = $newQueryC
$ws.Range("C4").Value2 = $newQueryC

# Update the sheet view: active selection moves from D4 to C3, with the
# window scrolled down so row 3 is at the top of the visible area.
$ws.Range("C3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1

$wb.Saved = $false
